# Save outputs to subfolder data/results/<model_name>/<val> or <test>
#
# This script:
#  1. Removes the "Papers" worksheet (its few rows are being folded back
#     into "Experiments" as a new "From Papers" section).
#  2. Bolds a handful of already-existing result cells.
#  3. Adds a new results row (row 5) for the histogram-rescaling ablation.
#  4. Adds a new "From Papers" section (rows 8-11) with a header row and
#     the Eigen et al. / DORN / Laina et al. reference numbers.
#  5. Leaves the active selection on A6, matching the new layout.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Drop the old "Papers" sheet -----------------------------------
[void]$wb.Worksheets("Papers").Delete()

$ws = $wb.Worksheets("Experiments")

# --- 2. Bold a few existing cells in row 3 -----------------------------
$ws.Range("D3").Font.Bold = $true
$ws.Range("E3").Font.Bold = $true
$ws.Range("H3").Font.Bold = $true

# --- Populate the brand-new label strings first, in the order they
#     should land in the shared-string table: "-", "Laina et. al.",
#     "?", "q", the histogram-rescaling label, then "From Papers".
$ws.Range("E10").Value = "-"
$ws.Range("A11").Value = "Laina et. al."
$ws.Range("E11").Value = "?"
$ws.Range("J10").Value = "q"
$ws.Range("A5").Value = "DORN + Histogram Rescaling (No noise, no albedo, no falloff)"
$ws.Range("A8").Value = "From Papers"

# --- 3. New row 5: DORN + Histogram Rescaling (no noise/albedo/falloff)
$ws.Range("B5").Value = 0.89902345485840895
$ws.Range("C5").Value = 0.96994871391908999
$ws.Range("D5").Value = 0.98954798995482196
$ws.Range("E5").Value = 0.3224211819335
$ws.Range("F5").Value = 0.48643683646165198
$ws.Range("G5").Value = 0.0947125232712816
$ws.Range("H5").Value = 0.088471443711215
$ws.Range("I5").Value = 0.0412984580063847
$ws.Range("J5").Value = 0.14785928237148899

$ws.Range("B5").Font.Bold = $true
$ws.Range("C5").Font.Bold = $true
$ws.Range("F5").Font.Bold = $true
$ws.Range("G5").Font.Bold = $true
$ws.Range("I5").Font.Bold = $true
$ws.Range("J5").Font.Bold = $true

# --- 4. New "From Papers" block (rows 8-11) ----------------------------
# Row 8: header row for the "From Papers" block
$ws.Range("B8").Value = "delta1"
$ws.Range("C8").Value = "delta2"
$ws.Range("D8").Value = "delta3"
$ws.Range("E8").Value = "mse"
$ws.Range("F8").Value = "rmse"
$ws.Range("G8").Value = "rel_abs_diff"
$ws.Range("H8").Value = "rel_sqr_diff"
$ws.Range("I8").Value = "log10"
$ws.Range("J8").Value = "log_rmse"

# Row 9: Eigen et al. (name only, numbers unknown)
$ws.Range("A9").Value = "Eigen et. al."

# Row 10: DORN, copied in from the old "Papers" sheet
$ws.Range("A10").Value = "DORN"
$ws.Range("B10").Value = 0.82799999999999996
$ws.Range("C10").Value = 0.96499999999999997
$ws.Range("D10").Value = 0.99199999999999999
$ws.Range("F10").Value = 0.50900000000000001
$ws.Range("G10").Value = 0.115
$ws.Range("H10").Value = "-"
$ws.Range("I10").Value = 0.050999999999999997

# Row 11: Laina et al.
$ws.Range("B11").Value = 0.81100000000000005
$ws.Range("C11").Value = 0.95299999999999996
$ws.Range("D11").Value = 0.98799999999999999
$ws.Range("F11").Value = 0.57299999999999995
$ws.Range("G11").Value = 0.127
$ws.Range("H11").Value = "?"
$ws.Range("I11").Value = 0.055
$ws.Range("J11").Value = "?"

# --- 5. Selection -------------------------------------------------------
[void]$ws.Range("A6").Select()
